$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description for "pop_region" (row 6)
$ws.Range("G6").Value = "Données démographiques et géographiques par région"

# Update the description for "temperature_moy" (row 5)
$ws.Range("G5").Value = "Évolution des températures mensuelles"

# Update the description for "immobilier_prix" (row 17)
$ws.Range("G17").Value = "Analyse des prix des biens immobiliers dans les grandes villes"

# Flip "transport_pub" (row 11) from open_data to closed_data
$ws.Range("F11").Value = "closed_data"

# Move the active selection, matching the saved workbook state
$ws.Range("D24").Select()
